$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the product descriptions in A2:A4 for the new product list.
$ws.Range("A2").Value = "Ar Condicionado Split Hw On/off Eco Garden Gree 18000 Btus, Quente/Frio, 220V, Monofásico – GWH18QD-D3NNB4B"
$ws.Range("A3").Value = "Ventilador de Parede 1 Metro, Com 3 velocidades, Ventisol, 220V"
$ws.Range("A4").Value = "Frigobar Midea MRC06B1, 45L, 110V, Branco"

# The old A4 text wrapped onto two lines, giving it a taller row; the new
# text fits on a single line, so let the row shrink back to the default height.
$ws.Rows.Item(4).AutoFit()

# Move the active selection down to A16, as in the edited workbook.
$ws.Range("A16").Select()
